$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain text formatting (values look numeric/percent)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "41.680.97"
$ws.Range("E2").Value = "  -0.15%  "

$ws.Range("D3").Value = "2.168.99"
$ws.Range("E3").Value = "  -1.95%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "238.06"
$ws.Range("E5").Value = "  -1.14%  "

$ws.Range("D6").Value = "0.607"
$ws.Range("E6").Value = "  -2.67%  "

$ws.Range("D7").Value = "71.31"
$ws.Range("E7").Value = "  -1.48%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "0.575"
$ws.Range("E9").Value = "  -2.98%  "

$ws.Range("D10").Value = "39.88"
$ws.Range("E10").Value = "  -4.35%  "

$ws.Range("D11").Value = "0.0904"
$ws.Range("E11").Value = "  -4.19%  "

$ws.Range("D12").Value = "54.23"
$ws.Range("E12").Value = "  -4.24%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.100"
$ws.Range("E13").Value = "  -3.06%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "6.70"
$ws.Range("E14").Value = "  -3.11%  "

$ws.Range("D15").Value = "2.496.87"
$ws.Range("E15").Value = "  -1.82%  "

$ws.Range("D16").Value = "14.28"
$ws.Range("E16").Value = "  +0.80%  "

$ws.Range("D17").Value = "2.165.62"
$ws.Range("E17").Value = "  -3.54%  "

$ws.Range("D18").Value = "0.785"
$ws.Range("E18").Value = "  -5.59%  "

$ws.Range("D19").Value = "41.546.30"
$ws.Range("E19").Value = "  -0.22%  "

$ws.Range("E20").Value = "  -2.98%  "

$ws.Range("D21").Value = "69.66"
$ws.Range("E21").Value = "  -3.37%  "

$ws.Range("D22").Value = "5.76"
$ws.Range("E22").Value = "  -6.08%  "

$ws.Range("D23").Value = "10.21"
$ws.Range("E23").Value = "  -4.62%  "

$ws.Range("D24").Value = "226.68"

$ws.Range("D25").Value = "1.97"
$ws.Range("E25").Value = "  -3.37%  "

$ws.Range("E26").Value = "  +0.25%  "

$ws.Range("D27").Value = "10.69"
$ws.Range("E27").Value = "  -5.48%  "

$ws.Range("D28").Value = "3.30"
$ws.Range("E28").Value = "  -8.66%  "

$ws.Range("D29").Value = "2.17"
$ws.Range("E29").Value = "  -4.08%  "

$ws.Range("E30").Value = "  -1.13%  "

$ws.Range("D31").Value = "170.84"
$ws.Range("E31").Value = "  +2.28%  "

$ws.Range("D32").Value = "19.78"
$ws.Range("E32").Value = "  -2.73%  "

$ws.Range("D33").Value = "32.67"
$ws.Range("E33").Value = "  +9.13%  "

$ws.Range("D34").Value = "0.0772"
$ws.Range("E34").Value = "  -2.20%  "

$ws.Range("D35").Value = "5.10"
$ws.Range("E35").Value = "  -7.16%  "

$ws.Range("D36").Value = "0.120"
$ws.Range("E36").Value = "  -3.07%  "

$ws.Range("E37").Value = "  -0.88%  "

$ws.Range("D38").Value = "4.30"
$ws.Range("E38").Value = "  +2.23%  "

$ws.Range("E39").Value = "  +1.16%  "

$ws.Range("D40").Value = "12.31"
$ws.Range("E40").Value = "  -7.98%  "

$ws.Range("D41").Value = "2.07"
$ws.Range("E41").Value = "  -2.00%  "

$ws.Range("D42").Value = "5.34"
$ws.Range("E42").Value = "  -4.76%  "

$ws.Range("D43").Value = "58.97"
$ws.Range("E43").Value = "  -7.46%  "

$ws.Range("D44").Value = "0.189"
$ws.Range("E44").Value = "  -3.16%  "

$ws.Range("D45").Value = "8.44"
$ws.Range("E45").Value = "  -2.57%  "

$ws.Range("D46").Value = "0.0966"
$ws.Range("E46").Value = "  -3.12%  "

$ws.Range("D47").Value = "97.14"
$ws.Range("E47").Value = "  -5.24%  "

$ws.Range("E48").Value = "  -2.04%  "

$ws.Range("D49").Value = "1.11"
$ws.Range("E49").Value = "  -4.04%  "

$ws.Range("D50").Value = "2.16"
$ws.Range("E50").Value = "  -6.63%  "

$ws.Range("E51").Value = "  -2.30%  "
